$d = $word.ActiveDocument

# --- Change 1: append a new sentence (as its own run) to the Find/Replace
# dialog paragraph, describing the new Find/Replace toolbar. ---
$rng1 = $d.Content
$rng1.Find.Execute(
    "find or replace dialog will reset the status bar contents. ",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng1.Collapse(0)
$rng1.InsertAfter("Additionally, a find and replace toolbar is available that sits at the bottom of the main application window. This toolbar can be accessed from the Find/Replace Bar menu option. ")

# --- Change 2: add a new bulleted "Known Issues" entry about the status bar
# disappearing when using the Find/Replace toolbar. ---
$rng2 = $d.Content
$rng2.Find.Execute(
    "The feature checking for unsaved changes is currently a little buggy sometimes. ",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng2.Collapse(0)
$rng2.InsertParagraphAfter()

$paras = $d.Paragraphs
$newPara = $paras.Item($paras.Count - 1)
$newRange = $newPara.Range
$newRange.InsertAfter("Status bar at bottom of window main sometimes temporarily disappear when using the Find/Replace toolbar. ")
